$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting D:K -> E:L
$ws.Columns("D").Insert()

# Copy number formats/styles from column E (the old D, now shifted) into the new column D
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the newly reported period values
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 4922900
$ws.Range("D9").Value = 2870600
$ws.Range("D10").Value = 2052300
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 28700
$ws.Range("D15").Value = 680500
$ws.Range("D17").Value = 4090700
$ws.Range("D18").Value = 832200
$ws.Range("D20").Value = 2800
$ws.Range("D21").Value = 1515600
$ws.Range("D22").Value = 127900
$ws.Range("D23").Value = 707100
$ws.Range("D24").Value = 159700
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 547400
$ws.Range("D27").Value = 547100
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = -300
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -2800
$ws.Range("D33").Value = 546900
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 546900
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 319300
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 609500
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 164100
$ws.Range("D46").Value = 1092900
$ws.Range("D47").Value = "NA"
$ws.Range("D48").Value = 5169000
$ws.Range("D49").Value = 6160300
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 205100
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 12627300
$ws.Range("D57").Value = 360000
$ws.Range("D58").Value = 20300
$ws.Range("D59").Value = 480400
$ws.Range("D60").Value = 860700
$ws.Range("D61").Value = 4153500
$ws.Range("D62").Value = 1153000
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 6172700
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 2264500
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 6454600
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 546900
$ws.Range("D83").Value = 680500
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 1411400
$ws.Range("D91").Value = -546100
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -1371800
$ws.Range("D96").Value = -152600
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -187600
$ws.Range("D101").Value = -1300
$ws.Range("D102").Value = -149300

# Row 94 (Capital Expenditures) also received a corrected/newly reported value in column E
$ws.Range("E94").Value = -860900